$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 174006
$ws.Range("C4").Value = 164000
$ws.Range("C5").Value = 10006
$ws.Range("C8").Value = 64.39
